$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'27.100.38"
$ws.Range("E2").Value = "'  -0.66%  "

$ws.Range("D3").Value = "'1.631.85"
$ws.Range("E3").Value = "'  -1.17%  "

$ws.Range("E4").Value = "'  +0.37%  "

$ws.Range("D5").Value = "'216.52"
$ws.Range("E5").Value = "'  -1.03%  "

$ws.Range("E6").Value = "'  +1.52%  "

$ws.Range("E7").Value = "'  +0.36%  "

$ws.Range("D8").Value = "'0.253"
$ws.Range("E8").Value = "'  -1.38%  "

$ws.Range("D9").Value = "'0.0624"
$ws.Range("E9").Value = "'  -0.58%  "

$ws.Range("D10").Value = "'20.07"
$ws.Range("E10").Value = "'  -0.86%  "

$ws.Range("D11").Value = "'0.0847"
$ws.Range("E11").Value = "'  -0.10%  "

$ws.Range("D12").Value = "'1.865.98"
$ws.Range("E12").Value = "'  -0.87%  "

$ws.Range("D13").Value = "'1.629.23"
$ws.Range("E13").Value = "'  -1.24%  "

$ws.Range("D14").Value = "'4.12"
$ws.Range("E14").Value = "'  -0.68%  "

$ws.Range("D15").Value = "'0.543"
$ws.Range("E15").Value = "'  +0.21%  "

$ws.Range("D16").Value = "'65.93"
$ws.Range("E16").Value = "'  -2.85%  "

$ws.Range("D17").Value = "'27.135.38"
$ws.Range("E17").Value = "'  -0.53%  "

$ws.Range("D18").Value = "'0.0₃0734"
$ws.Range("E18").Value = "'  -0.70%  "

$ws.Range("D19").Value = "'214.76"
$ws.Range("E19").Value = "'  -3.32%  "

$ws.Range("D21").Value = "'6.84"
$ws.Range("E21").Value = "'  +0.95%  "

$ws.Range("D22").Value = "'4.41"
$ws.Range("E22").Value = "'  -0.91%  "

$ws.Range("D23").Value = "'2.51"
$ws.Range("E23").Value = "'  +0.89%  "

$ws.Range("D24").Value = "'9.10"
$ws.Range("E24").Value = "'  -1.69%  "

$ws.Range("D25").Value = "'147.47"
$ws.Range("E25").Value = "'  -0.17%  "

$ws.Range("E26").Value = "'  +0.37%  "

$ws.Range("D27").Value = "'7.37"
$ws.Range("E27").Value = "'  -0.49%  "

$ws.Range("E28").Value = "'  -1.60%  "

$ws.Range("D29").Value = "'15.61"
$ws.Range("E29").Value = "'  -1.74%  "

$ws.Range("E30").Value = "'  -0.09%  "

$ws.Range("E31").Value = "'  -0.70%  "

$ws.Range("D32").Value = "'3.36"
$ws.Range("E32").Value = "'  -0.02%  "

$ws.Range("D33").Value = "'3.02"
$ws.Range("E33").Value = "'  -0.69%  "

$ws.Range("D34").Value = "'1.306.62"
$ws.Range("E34").Value = "'  +2.55%  "

$ws.Range("D35").Value = "'1.56"
$ws.Range("E35").Value = "'  -1.81%  "

$ws.Range("E36").Value = "'  +0.69%  "

$ws.Range("D37").Value = "'0.0176"
$ws.Range("E37").Value = "'  -0.80%  "

$ws.Range("D38").Value = "'0.545"
$ws.Range("E38").Value = "'  +0.20%  "

$ws.Range("D39").Value = "'0.846"
$ws.Range("E39").Value = "'  +0.12%  "

$ws.Range("E40").Value = "'  +0.34%  "

$ws.Range("E41").Value = "'  +3.92%  "

$ws.Range("D42").Value = "'0.807"
$ws.Range("E42").Value = "'  -0.39%  "

$ws.Range("D43").Value = "'5.33"
$ws.Range("E43").Value = "'  -1.09%  "

$ws.Range("D44").Value = "'1.776.41"
$ws.Range("E44").Value = "'  -0.87%  "

$ws.Range("D45").Value = "'61.95"
$ws.Range("E45").Value = "'  -2.21%  "

$ws.Range("D46").Value = "'90.67"
$ws.Range("E46").Value = "'  -2.04%  "

$ws.Range("E47").Value = "'  -0.32%  "

$ws.Range("D48").Value = "'0.0₆0100"
$ws.Range("E48").Value = "'  -1.98%  "

$ws.Range("D49").Value = "'0.0513"
$ws.Range("E49").Value = "'  -0.14%  "

$ws.Range("B50").Value = "'WEMIXToken"
$ws.Range("C50").Value = "'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D50").Value = "'0.763"
$ws.Range("E50").Value = "'  +14.45%  "

$ws.Range("B51").Value = "'EnergySwap"
$ws.Range("C51").Value = "'https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D51").Value = "'7.61"
$ws.Range("E51").Value = "'  -1.66%  "

